# Insert a new "Industry" column between "Stock Name" (B) and "Mutual Fund" (C),
# shifting Mutual Fund, Status, Jan_2026, Dec_2025, Oct_2025, MoM, QoQ one column right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at C - this shifts existing C..I to D..J
$ws.Columns.Item(3).Insert()

# Header for the new column
$ws.Cells.Item(1, 3).Value = "Industry"

$industryByRow = @{
    2 = "Banks"
    3 = "Banks"
    4 = "Capital Markets"
    5 = "Banks"
    6 = "Banks"
    7 = "Finance"
    8 = "Capital Markets"
    9 = "Finance"
    10 = "Capital Markets"
    11 = "Banks"
    12 = "Finance"
    13 = "Banks"
    14 = "Finance"
    15 = "Capital Markets"
    16 = "Financial Technology (Fintech)"
    17 = "Capital Markets"
    18 = "Banks"
    19 = "Insurance"
    20 = "Banks"
    21 = "Banks"
    22 = "Banks"
    23 = "Insurance"
    24 = "Financial Technology (Fintech)"
    25 = "Capital Markets"
    26 = "Finance"
    27 = "Finance"
    28 = "Insurance"
    29 = "Finance"
    30 = "Insurance"
    31 = "Capital Markets"
    32 = "Capital Markets"
    33 = "Finance"
    34 = "Capital Markets"
    35 = "Capital Markets"
    36 = "Banks"
}

foreach ($row in $industryByRow.Keys) {
    $ws.Cells.Item($row, 3).Value = $industryByRow[$row]
}
